# Auto-generated edit script: updates the 100 arithmetic expressions
# in the single 20x5 table, in row-major (document) order, to match
# the target diff. Each table cell holds exactly one run of text, so
# setting Cell.Range.Text preserves the existing run/paragraph formatting.

$d = $word.ActiveDocument

$newValues = @("15-11=","43-27=","61-10=","71-7=","17+32=","97-7=","66-52=","53-8=","1+47=","85-59=","63-63=","32+20=","30+27=","8-8=","97-32=","77-70=","31+24=","16-2=","56-14=","69+0=","24+36=","85-24=","73-57=","85-47=","46+17=","53-7=","9+34=","60+5=","82-18=","99-40=","89-7=","14+63=","26+68=","97-77=","87-20=","22+23=","45+29=","10+15=","67-20=","8+24=","35-24=","67+18=","29+35=","77+12=","13+61=","65-51=","66-61=","45+23=","68-32=","85-83=","15+20=","70-22=","94-52=","65-39=","58+28=","60+20=","42-25=","78-43=","95-78=","74-39=","4+4=","96-9=","91-23=","67+1=","64+29=","24-2=","8+64=","48-22=","53-24=","27+71=","98-74=","74-14=","51+45=","74-43=","94-65=","79-4=","6+43=","1+49=","65-45=","55+12=","4+3=","64+1=","79-24=","3+63=","13+24=","57-40=","30+58=","77-29=","78-47=","35+11=","81-77=","76-18=","70-17=","75-66=","90-21=","41+28=","41+17=","69-50=","92-91=","34+23=")

$table = $d.Tables.Item(1)
$rows = $table.Rows.Count
$cols = $table.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    $row = $table.Rows.Item($r)
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $row.Cells.Item($c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Output "Updated $i cells"
